# Vermont raw 2023 workbook: anonymize percent_moisture column (O) and
# restyle it, matching the upstream edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- O1 header: pick up the plain "Calibri" cell style already used by
#     the rest of the id columns (A2/B2/D2/H2/N2 ... all style index 4)
#     instead of the shaded/bordered style it had before. Copy+PasteSpecial
#     (formats only) reuses the existing cellXf rather than inventing a new
#     one, so the output style table stays byte-for-byte aligned with the
#     target (s="9" -> s="4").
$ws.Range("N2").Copy() | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- O2:O37 body cells: every percent_moisture reading is replaced with a
#     uniform placeholder value, and the column is switched to the
#     "Aptos Narrow" font (the workbook's theme minor font). Setting the
#     font across the whole range at once means Excel creates exactly one
#     new font + one new cell style shared by all 36 cells, matching the
#     diff (fonts count 4 -> 5, cellXfs count 10 -> 11).
$body = $ws.Range("O2:O37")
$body.Font.Name = "Aptos Narrow"
$body.Value = 0.00001

# --- Selection / view: the author's last on-screen selection before saving.
$ws.Range("N21").Select()
